$d = $word.ActiveDocument

# --- 1. Remove the "_GoBack" bookmark that currently sits in the middle
#        of the "...and displaying [bookmark]that data is what we'll
#        cover in the next lesson." paragraph. Word still exposes this
#        hidden bookmark by name even though it is excluded from
#        Bookmarks.Count. Deleting it only removes the bookmarkStart/
#        bookmarkEnd markers and leaves every run/word around it
#        untouched, so the two sentences stay in the same paragraph.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2. Locate the "So I'll see you there." paragraph, then step to the
#        very next paragraph, which is the one holding the legacy OLE /
#        ActiveX control object (the little form-control icon at the
#        end of the document). Removing that paragraph's Range (which
#        includes its own trailing paragraph mark) deletes the control
#        and merges the paragraph away, leaving the already-empty final
#        paragraph right after it untouched.
$seeYou = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd([char]13) -eq "So I'll see you there.") {
        $seeYou = $candidate
    }
}

if ($seeYou -ne $null) {
    $olePara = $seeYou.Next()
    $olePara.Range.Delete()
}

# --- 3. Re-create the "_GoBack" bookmark, collapsed (empty), inside the
#        now-trailing empty paragraph (the document's original, already
#        empty final paragraph) so the document keeps exactly one
#        "_GoBack" bookmark, just relocated to the very end.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(1)
$d.Bookmarks.Add("_GoBack", $tail)
